# Auto update Excel log
# Appends newly-logged sensor readings to the PIR, Humidity and
# Temperature sheets (SeniorConnect_MasterLog.xlsx).
#
# All of these sheets store every column as literal text (dates like
# "2026-02-01", times like "18:26:33", percentages like "78.4%" and
# temperatures like "29.6C" are plain strings, not real Excel
# number/date/percentage values) so each new row's range is formatted
# as Text ("@") before the values are written in order to stop Excel's
# automatic type-inference from turning them into dates/numbers.

$wb = $excel.ActiveWorkbook

# NOTE: named parameters don't bind reliably in this PowerShell host, so
# this helper is called positionally: Append-LogRows <sheet> <startRow> <rows>
function Append-LogRows($SheetName, $StartRow, $Rows) {
    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow
    foreach ($row in $Rows) {
        $rng = $ws.Range("A$r`:F$r")
        $rng.NumberFormat = "@"
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }
}

# ---- PIR sheet: rows 41-45 (Date, Timestamp, Hour, Location, Value, Status) ----
$pirRows = @(
    ,@("2026-02-01", "18:26:33", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-02-01", "18:26:35", "18:00", "Bathroom", "Motion Detected", "Active")
    ,@("2026-02-01", "18:26:35", "18:00", "Bathroom", "No Motion", "Inactive")
    ,@("2026-02-01", "18:26:35", "18:00", "Bathroom", "Motion Detected", "Active")
    ,@("2026-02-01", "18:27:29", "18:00", "Bathroom", "Motion Detected", "Active")
)
Append-LogRows "PIR" 41 $pirRows

# ---- Humidity sheet: rows 80-93 ----
$humidityRows = @(
    ,@("2026-02-01", "18:26:30", "18:00", "Bathroom", "78.4%", "Active")
    ,@("2026-02-01", "18:26:31", "18:00", "Bathroom", "79.4%", "Active")
    ,@("2026-02-01", "18:26:31", "18:00", "Bathroom", "79.4%", "Active")
    ,@("2026-02-01", "18:26:32", "18:00", "Bathroom", "80.4%", "Active")
    ,@("2026-02-01", "18:26:36", "18:00", "Bathroom", "81.4%", "Active")
    ,@("2026-02-01", "18:26:37", "18:00", "Bathroom", "80.6%", "Active")
    ,@("2026-02-01", "18:26:43", "18:00", "Bathroom", "78.2%", "Active")
    ,@("2026-02-01", "18:26:49", "18:00", "Bathroom", "78.9%", "Active")
    ,@("2026-02-01", "18:26:53", "18:00", "Bathroom", "77.9%", "Active")
    ,@("2026-02-01", "18:27:08", "18:00", "Bathroom", "85.4%", "Active")
    ,@("2026-02-01", "18:27:13", "18:00", "Bathroom", "92.2%", "Active")
    ,@("2026-02-01", "18:27:18", "18:00", "Bathroom", "93.2%", "Active")
    ,@("2026-02-01", "18:27:23", "18:00", "Bathroom", "90.8%", "Active")
    ,@("2026-02-01", "18:27:28", "18:00", "Bathroom", "88.5%", "Active")
)
Append-LogRows "Humidity" 80 $humidityRows

# ---- Temperature sheet: rows 80-93 ----
$temperatureRows = @(
    ,@("2026-02-01", "18:26:30", "18:00", "Bathroom", "29.6C", "Active")
    ,@("2026-02-01", "18:26:31", "18:00", "Bathroom", "29.6C", "Active")
    ,@("2026-02-01", "18:26:32", "18:00", "Bathroom", "29.6C", "Active")
    ,@("2026-02-01", "18:26:33", "18:00", "Bathroom", "29.6C", "Active")
    ,@("2026-02-01", "18:26:36", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:26:37", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:26:44", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:26:49", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:26:54", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:27:09", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:27:14", "18:00", "Bathroom", "29.7C", "Active")
    ,@("2026-02-01", "18:27:19", "18:00", "Bathroom", "29.8C", "Active")
    ,@("2026-02-01", "18:27:24", "18:00", "Bathroom", "29.9C", "Active")
    ,@("2026-02-01", "18:27:29", "18:00", "Bathroom", "29.9C", "Active")
)
Append-LogRows "Temperature" 80 $temperatureRows
